$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 36 and 37 (columns F:V) ---
$r36 = $ws.Range("F36:V36")
$r37 = $ws.Range("F37:V37")
$v36 = $r36.Value2
$v37 = $r37.Value2
$r36.Value = $v37
$r37.Value = $v36

# --- Swap rows 76 and 77 (columns F:V) ---
$r76 = $ws.Range("F76:V76")
$r77 = $ws.Range("F77:V77")
$v76 = $r76.Value2
$v77 = $r77.Value2
$r76.Value = $v77
$r77.Value = $v76

# --- Add new row 78 with formatting copied from row 77 ---
$ws.Range("A77:V77").Copy()
$ws.Range("A78:V78").PasteSpecial(-4122)

$ws.Range("A78").Value = 77
$ws.Range("B78").Value = "bosnia-and-herzegovina"
$ws.Range("C78").Value = "premijer-liga-bih"
$ws.Range("D78").Value = "2023-2024"
$ws.Range("E78").Value = 45235.82291666666
$ws.Range("F78").Value = "FK Sarajevo"
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = "Borac Banja Luka"
$ws.Range("I78").Value = 1
$ws.Range("J78").Value = 1.78
$ws.Range("K78").Value = "04/11/2023 08:12"
$ws.Range("L78").Value = 2.2
$ws.Range("M78").Value = "05/11/2023 19:18"
$ws.Range("N78").Value = 3.27
$ws.Range("O78").Value = "04/11/2023 08:12"
$ws.Range("P78").Value = 3.34
$ws.Range("Q78").Value = "05/11/2023 19:18"
$ws.Range("R78").Value = 4.05
$ws.Range("S78").Value = "04/11/2023 08:12"
$ws.Range("T78").Value = 3.15
$ws.Range("U78").Value = "05/11/2023 19:18"
$ws.Range("V78").Value = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/fk-sarajevo-borac-banja-luka/nNgZ3hOH/"

Write-Output "done"
